$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: Summary
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = "Noura Al Mazroui"
$summary.Range("B4").Value = 1627.28
$summary.Range("B6").Value = 725773
$summary.Range("B7").Value = 500446
$summary.Range("B8").Value = 225327
$summary.Range("B9").Value = 1.45

# ---------------------------------------------------------------------
# Sheet 2: Assets
# ---------------------------------------------------------------------
$assets = $wb.Worksheets.Item("Assets")

# Update existing vehicle row (Mid-range Car -> Premium Car)
$assets.Range("B2").Value = "Premium Car"
$assets.Range("C2").Value = 254129

# Insert a new row for the second vehicle (shifts Liquid Assets + Total down)
$assets.Rows.Item(3).Insert()

# Copy formatting from row 2 (Vehicles data row) into the new row 3
$assets.Range("A2:C2").Copy($assets.Range("A3:C3"))

$assets.Range("A3").Value = "Vehicles"
$assets.Range("B3").Value = "Luxury Car"
$assets.Range("C3").Value = 468951

# Liquid assets row (now row 4)
$assets.Range("C4").Value = 2693

# Total assets row (now row 5)
$assets.Range("C5").Value = 725773

# ---------------------------------------------------------------------
# Sheet 3: Liabilities
# ---------------------------------------------------------------------
$liab = $wb.Worksheets.Item("Liabilities")

# Insert three new rows before the existing Credit Cards row (row 2)
$liab.Rows.Item(2).Resize(3).Insert()

# Copy formatting from the (now shifted) data row 5 into the new rows 2-4
$liab.Range("A5:E5").Copy($liab.Range("A2:E2"))
$liab.Range("A5:E5").Copy($liab.Range("A3:E3"))
$liab.Range("A5:E5").Copy($liab.Range("A4:E4"))

# Row 2: Auto Loans / Vehicle Loan 1
$liab.Range("A2").Value = "Auto Loans"
$liab.Range("B2").Value = "Vehicle Loan 1"
$liab.Range("C2").Value = 152477
$liab.Range("D2").Value = 2118
$liab.Range("E2").Value = 6

# Row 3: Auto Loans / Vehicle Loan 2
$liab.Range("A3").Value = "Auto Loans"
$liab.Range("B3").Value = "Vehicle Loan 2"
$liab.Range("C3").Value = 281371
$liab.Range("D3").Value = 3350
$liab.Range("E3").Value = 7

# Row 4: Personal Loans / Personal Loan
$liab.Range("A4").Value = "Personal Loans"
$liab.Range("B4").Value = "Personal Loan"
$liab.Range("C4").Value = 58558
$liab.Range("D4").Value = 1220
$liab.Range("E4").Value = 4

# Row 5: Credit Cards / Credit Card Balance (existing row, now shifted down)
$liab.Range("C5").Value = 8040
$liab.Range("D5").Value = 402

# Row 6: TOTAL LIABILITIES (existing row, now shifted down)
$liab.Range("C6").Value = 500446
